$d = $word.ActiveDocument

# --- Paragraph: "Atualmente composta por uma unica pessoa ..." -----------
# 1) Fix the missing space between sentences: "presencialmente.A" -> "presencialmente. A"
$d.Content.Find.Execute(
    "presencialmente.A empresa", $false, $false, $false, $false, $false,
    $true, 1, $false, "presencialmente. A empresa", 2) | Out-Null

# 2) Italicize the word "WhatsApp" in that same paragraph.
$rng = $d.Content
$rng.Find.Execute("WhatsApp") | Out-Null
$rng.Italic = 1

# 3) Italicize the supplier list "Ruby Rose, ..., Pink 21".
$rng2 = $d.Content
$rng2.Find.Execute(
    "Ruby Rose, Fenzza, Di Grezzo, Miss Lary, Dalla, Luisance, Mahav, SP Colors, Macrilan, Pink 21") | Out-Null
$rng2.Italic = 1

# --- Paragraph: "O projeto de desenvolvimento e implementacao do Loculus System ..." ---
# 4) Italicize the second occurrence of "Loculus System" (the first one is the title).
$rng3 = $d.Content
$rng3.Find.Execute("Loculus System") | Out-Null
$rng3.Find.Execute("Loculus System") | Out-Null
$rng3.Italic = 1
